# -----------------------------------------------------------------------
# Commit: "new sim results and new calculation"
#
# 1) Insert a new worksheet "sharpe_period" immediately before "VaR"
#    (so the final tab order is: annualised_return, mean_period_return,
#     sharpe_annualized, sharpe_period, VaR) and populate it with the
#     new correlation/flag matrix.
# 2) Overwrite a batch of cell values on the four pre-existing sheets
#    (annualised_return, mean_period_return, sharpe_annualized, VaR)
#    with the new simulation results.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "sharpe_period" sheet right before "VaR"
# ---------------------------------------------------------------------
$varSheet = $wb.Worksheets.Item("VaR")
$newSheet = $wb.Worksheets.Add($varSheet)
$newSheet.Name = "sharpe_period"

$labels = @("minvar_ports_equalw","minvar_w_cryptos_ahc_equalw","minvar_w_cryptos_kmeans_equalw","minvar_w_cryptos_kshape_equalw","minvar_w_cryptos_random_equalw","rand_ports_equalw","random_w_cryptos_ahc_equalw","random_w_cryptos_kmeans_equalw","random_w_cryptos_kshape_equalw","random_w_cryptos_random_equalw")
$cols = @("B","C","D","E","F","G","H","I","J","K")

# header row (B1:K1) + row labels (A2:A11)
for ($i = 0; $i -lt 10; $i++) {
    $newSheet.Range($cols[$i] + "1").Value = $labels[$i]
    $newSheet.Range("A" + ($i + 2)).Value = $labels[$i]
}

$sharpePeriodData = @{
    "B2"="1";  "C2"="0";       "D2"="0";      "E2"="0";      "F2"="0";  "G2"="0";      "H2"="0"; "I2"="0";      "J2"="0";    "K2"="0";
    "B3"="0";  "C3"="1";       "D3"="0";      "E3"="0";      "F3"="1";  "G3"="0";      "H3"="0"; "I3"="0";      "J3"="0";    "K3"="0";
    "B4"="0";  "C4"="0";       "D4"="1";      "E4"="1";      "F4"="0";  "G4"="0.0001"; "H4"="0"; "I4"="0";      "J4"="0";    "K4"="0";
    "B5"="0";  "C5"="0";       "D5"="1";      "E5"="1";      "F5"="0";  "G5"="0.0001"; "H5"="0"; "I5"="0";      "J5"="0";    "K5"="0";
    "B6"="0";  "C6"="1";       "D6"="0";      "E6"="0";      "F6"="1";  "G6"="0";      "H6"="0"; "I6"="0";      "J6"="0";    "K6"="0";
    "B7"="0";  "C7"="0";       "D7"="0.0001"; "E7"="0.0001"; "F7"="0";  "G7"="1";      "H7"="0"; "I7"="0";      "J7"="0";    "K7"="0";
    "B8"="0";  "C8"="0";       "D8"="0";      "E8"="0";      "F8"="0";  "G8"="0";      "H8"="1"; "I8"="1";      "J8"="1";    "K8"="1";
    "B9"="0";  "C9"="0";       "D9"="0";      "E9"="0";      "F9"="0";  "G9"="0";      "H9"="1"; "I9"="1";      "J9"="1";    "K9"="0.4533";
    "B10"="0"; "C10"="0";      "D10"="0";     "E10"="0";     "F10"="0"; "G10"="0";     "H10"="1";"I10"="1";     "J10"="1";   "K10"="0.592";
    "B11"="0"; "C11"="0";      "D11"="0";     "E11"="0";     "F11"="0"; "G11"="0";     "H11"="1";"I11"="0.4533";"J11"="0.592";"K11"="1";
}
foreach ($ref in $sharpePeriodData.Keys) {
    $newSheet.Range($ref).Value = [double]$sharpePeriodData[$ref]
}

# ---------------------------------------------------------------------
# Step 2: new simulation values on the pre-existing sheets
# ---------------------------------------------------------------------
$edits = @{
    "annualised_return" = @{
        "C2"="0.002";  "F2"="0.1086"; "G2"="0";
        "B3"="0.002";  "D3"="0";      "E3"="0";      "F3"="1";
        "C4"="0";      "F4"="0";      "K4"="1";
        "C5"="0";      "F5"="0";      "K5"="1";
        "B6"="0.1086"; "C6"="1";      "D6"="0";      "E6"="0";
        "B7"="0";      "H7"="0.5229"; "J7"="0";
        "G8"="0.5229"; "I8"="0.302";  "J8"="0";      "K8"="0";
        "H9"="0.302";  "J9"="0.2281"; "K9"="0";
        "G10"="0";     "H10"="0";     "I10"="0.2281";"K10"="0";
        "D11"="1";     "E11"="1";     "H11"="0";     "I11"="0"; "J11"="0";
    }
    "mean_period_return" = @{
        "C2"="0.002";  "F2"="0.1086"; "G2"="0";
        "B3"="0.002";  "D3"="0";      "E3"="0";      "F3"="1";
        "C4"="0";      "F4"="0";      "K4"="1";
        "C5"="0";      "F5"="0";      "K5"="1";
        "B6"="0.1086"; "C6"="1";      "D6"="0";      "E6"="0";
        "B7"="0";      "H7"="0.5229"; "J7"="0";
        "G8"="0.5229"; "I8"="0.302";  "J8"="0";      "K8"="0";
        "H9"="0.302";  "J9"="0.2281"; "K9"="0";
        "G10"="0";     "H10"="0";     "I10"="0.2281";"K10"="0";
        "D11"="1";     "E11"="1";     "H11"="0";     "I11"="0"; "J11"="0";
    }
    "sharpe_annualized" = @{
        "D3"="0";  "E3"="0";  "F3"="1"; "K3"="1";
        "C4"="0";  "F4"="0";  "H4"="1"; "I4"="1"; "J4"="0.0015";
        "C5"="0";  "F5"="0";  "H5"="1"; "I5"="1"; "J5"="0.0016";
        "C6"="1";  "D6"="0";  "E6"="0"; "H6"="0"; "I6"="0"; "K6"="0.3249";
        "H7"="0";  "J7"="0";  "K7"="0";
        "D8"="1";  "E8"="1";  "F8"="0"; "G8"="0"; "J8"="0.0005999999999999999"; "K8"="0";
        "D9"="1";  "E9"="1";  "F9"="0"; "J9"="0.0321"; "K9"="0";
        "D10"="0.0015"; "E10"="0.0016"; "G10"="0"; "H10"="0.0005999999999999999"; "I10"="0.0321"; "K10"="0";
        "C11"="1"; "F11"="0.3249"; "G11"="0"; "H11"="0"; "I11"="0"; "J11"="0";
    }
    "VaR" = @{
        "K3"="1";
        "H4"="1"; "I4"="1"; "J4"="0.0062";
        "H5"="1"; "I5"="1"; "J5"="0.0065";
        "K6"="0.7436";
        "D8"="1"; "E8"="1"; "I8"="1"; "J8"="0.0008";
        "D9"="1"; "E9"="1"; "H9"="1"; "J9"="0.033";
        "D10"="0.0062"; "E10"="0.0065"; "H10"="0.0008"; "I10"="0.033";
        "C11"="1"; "F11"="0.7436";
    }
}

foreach ($sheetName in $edits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $sheetEdits = $edits[$sheetName]
    foreach ($ref in $sheetEdits.Keys) {
        $ws.Range($ref).Value = [double]$sheetEdits[$ref]
    }
}
